# This script updates computed market/profit figures (columns H-N) on several
# "Leve profit" worksheets (one worksheet per crafting class) to refresh them
# with newer Universalis market-board pricing data, per the scheduled data-
# refresh run. Only numeric result cells are touched; no rows/columns are
# inserted or removed and no other content is modified.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 1424.4445
$ws.Range("I6").Value = 205
$ws.Range("J6").Value = 2400
$ws.Range("K6").Value = 615
$ws.Range("L6").Value = 7200
$ws.Range("M6").Value = -503
$ws.Range("N6").Value = -7424
# Row 8: On the Drip / Eye Drops
$ws.Range("H8").Value = 1474.6154
$ws.Range("I8").Value = 400
$ws.Range("J8").Value = 1797
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 5391
$ws.Range("M8").Value = -1061
$ws.Range("N8").Value = -5669

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 12: Strait Ain't the Gate / Bronze Scutum
$ws.Range("H12").Value = 27332.334
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 40997
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 40997
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -41343
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1817.8214
$ws.Range("I45").Value = 1103.96
$ws.Range("K45").Value = 1103.96
$ws.Range("M45").Value = -726.96
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 2526.111
$ws.Range("I102").Value = 2123.125
$ws.Range("K102").Value = 2123.125
$ws.Range("M102").Value = -501.125
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1063.6757
$ws.Range("I110").Value = 513.8182
$ws.Range("J110").Value = 5600
$ws.Range("K110").Value = 513.8182
$ws.Range("L110").Value = 5600
$ws.Range("M110").Value = 1531.1818
$ws.Range("N110").Value = -9690
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3079.9473
$ws.Range("I122").Value = 2052.5
$ws.Range("J122").Value = 4221.5557
$ws.Range("K122").Value = 6157.5
$ws.Range("L122").Value = 12664.6671
$ws.Range("M122").Value = -3707.5
$ws.Range("N122").Value = -17564.6671

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 9: I, Gladius / Viking Sword
$ws.Range("H9").Value = 34999
$ws.Range("J9").Value = 34999
$ws.Range("L9").Value = 34999
$ws.Range("N9").Value = -35335
# Row 11: Down on the Pharm / Amateur's Mortar
$ws.Range("H11").Value = 4649.6
$ws.Range("I11").Value = 82.666664
$ws.Range("J11").Value = 11500
$ws.Range("K11").Value = 82.666664
$ws.Range("L11").Value = 11500
$ws.Range("M11").Value = 57.333336
$ws.Range("N11").Value = -11780
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1631.4783
$ws.Range("I105").Value = 1330
$ws.Range("K105").Value = 1330
$ws.Range("M105").Value = 417
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1957.8889
$ws.Range("I107").Value = 1017.4286
$ws.Range("J107").Value = 5249.5
$ws.Range("K107").Value = 1017.4286
$ws.Range("L107").Value = 5249.5
$ws.Range("M107").Value = 902.5714
$ws.Range("N107").Value = -9089.5

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 19: Shielding Sales / Square Ash Shield
$ws.Range("H19").Value = 337
$ws.Range("I19").Value = 136.66667
$ws.Range("J19").Value = 537.3333
$ws.Range("K19").Value = 136.66667
$ws.Range("L19").Value = 537.3333
$ws.Range("M19").Value = 33.33332999999999
$ws.Range("N19").Value = -877.3333
# Row 24: What You Need / Square Ash Shield
$ws.Range("H24").Value = 337
$ws.Range("I24").Value = 136.66667
$ws.Range("J24").Value = 537.3333
$ws.Range("K24").Value = 136.66667
$ws.Range("L24").Value = 537.3333
$ws.Range("M24").Value = 33.33332999999999
$ws.Range("N24").Value = -877.3333
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3221.513
$ws.Range("I31").Value = 2297.9312
$ws.Range("J31").Value = 5899.9
$ws.Range("K31").Value = 2297.9312
$ws.Range("L31").Value = 5899.9
$ws.Range("M31").Value = -2002.9312
$ws.Range("N31").Value = -6489.9
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3221.513
$ws.Range("I34").Value = 2297.9312
$ws.Range("J34").Value = 5899.9
$ws.Range("K34").Value = 2297.9312
$ws.Range("L34").Value = 5899.9
$ws.Range("M34").Value = -2095.9312
$ws.Range("N34").Value = -6303.9
# Row 63: So You Think You Can Lance? / Mythrite Trident
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372
# Row 66: Sticks and Stones (L) / Mythrite Trident
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864
# Row 68: Do You Even String Bow / Holy Cedar Composite Bow
$ws.Range("H68").Value = 32180
$ws.Range("I68").Value = 15000
$ws.Range("J68").Value = 36475
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 36475
$ws.Range("M68").Value = -14251
$ws.Range("N68").Value = -37973
# Row 71: Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws.Range("H71").Value = 32180
$ws.Range("I71").Value = 15000
$ws.Range("J71").Value = 36475
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 109425
$ws.Range("M71").Value = -41256
$ws.Range("N71").Value = -116913
# Row 75: The Darkest Hearth / Dark Chestnut Spinning Wheel
$ws.Range("H75").Value = 32000
$ws.Range("J75").Value = 32000
$ws.Range("L75").Value = 32000
$ws.Range("N75").Value = -33996
# Row 78: Fruit of the Loom (L) / Dark Chestnut Spinning Wheel
$ws.Range("H78").Value = 32000
$ws.Range("J78").Value = 32000
$ws.Range("L78").Value = 96000
$ws.Range("N78").Value = -105984
# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Range("H80").Value = 26000
$ws.Range("J80").Value = 26000
$ws.Range("L80").Value = 26000
$ws.Range("N80").Value = -28246
# Row 81: Don't Ask Wyvern / Hallowed Chestnut Composite Bow
$ws.Range("H81").Value = 56776
$ws.Range("J81").Value = 56776
$ws.Range("L81").Value = 56776
$ws.Range("N81").Value = -58772
# Row 82: Aim to Please / Hallowed Chestnut Mask of Aiming
$ws.Range("H82").Value = 21795.25
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 21795.25
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 21795.25
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -22517.25
# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Range("H83").Value = 26000
$ws.Range("J83").Value = 26000
$ws.Range("L83").Value = 78000
$ws.Range("N83").Value = -89232
# Row 84: A Sky Pirate's Life for Me (L) / Hallowed Chestnut Composite Bow
$ws.Range("H84").Value = 56776
$ws.Range("J84").Value = 56776
$ws.Range("L84").Value = 170328
$ws.Range("N84").Value = -180312
# Row 85: To Protect My City, I Must Wear a Mask (L) / Hallowed Chestnut Mask of Aiming
$ws.Range("H85").Value = 21795.25
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 21795.25
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 21795.25
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -24291.25
# Row 87: Anatomy of a Drill Bit / Dragonscale Grinding Wheel
$ws.Range("H87").Value = 31555.445
$ws.Range("J87").Value = 31555.445
$ws.Range("L87").Value = 31555.445
$ws.Range("N87").Value = -33927.445
# Row 90: Pulling Them to the Grind (L) / Dragonscale Grinding Wheel
$ws.Range("H90").Value = 31555.445
$ws.Range("J90").Value = 31555.445
$ws.Range("L90").Value = 94666.33499999999
$ws.Range("N90").Value = -106522.335

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2138.3333
$ws.Range("I102").Value = 1587.95
$ws.Range("K102").Value = 1587.95
$ws.Range("M102").Value = 34.04999999999995
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2899.8
$ws.Range("I113").Value = 1099.6666
$ws.Range("K113").Value = 1099.6666
$ws.Range("M113").Value = 1070.3334
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2355.2368
$ws.Range("I126").Value = 1375.96
$ws.Range("J126").Value = 4238.4614
$ws.Range("K126").Value = 4127.88
$ws.Range("L126").Value = 12715.3842
$ws.Range("M126").Value = -1657.88
$ws.Range("N126").Value = -17655.3842

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 19: Targe Up / Leather Targe
$ws.Range("H19").Value = 401.5
$ws.Range("I19").Value = 401.5
$ws.Range("K19").Value = 401.5
$ws.Range("M19").Value = -231.5
# Row 30: Packing a Punch / Goatskin Cesti
$ws.Range("H30").Value = 3000
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2892
$ws.Range("N30").ClearContents()

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1473.1875
$ws.Range("I113").Value = 658.875
$ws.Range("J113").Value = 2287.5
$ws.Range("K113").Value = 1976.625
$ws.Range("L113").Value = 6862.5
$ws.Range("M113").Value = 193.375
$ws.Range("N113").Value = -11202.5

Write-Host "Updated 194 cells across 7 worksheets."
